# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" / "Valor Mora" pairs in rows 17-51 were re-ordered:
# the periods used to run from the newest (2405) down to the oldest
# (2107); they now run from the oldest (2107) up to the newest (2405).
# Reverse the 35-row block (E17:F51) so the labels and their amounts
# travel together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 17
$lastRow = 51

$periods = @()
$amounts = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $periods += $ws.Cells.Item($r, 5).Value()
    $amounts += $ws.Cells.Item($r, 6).Value()
}

$count = $periods.Count
for ($i = 0; $i -lt $count; $i++) {
    $r = $firstRow + $i
    $srcIndex = $count - 1 - $i
    $ws.Cells.Item($r, 5).Value = $periods[$srcIndex]
    $ws.Cells.Item($r, 6).Value = $amounts[$srcIndex]
}
